# Auto-generated script applying value updates to Tonberry_Profits workbook
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 13989.083
$ws.Range("I18").Value = 1720
$ws.Range("K18").Value = 1720
$ws.Range("M18").Value = -1436
$ws.Range("H132").Value = 1181.2609
$ws.Range("I132").Value = 1181.2609
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3543.7827
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1013.7827
$ws.Range("N132").ClearContents()
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H135").Value = 776.8889
$ws.Range("I135").Value = 686.5
$ws.Range("K135").Value = 6178.5
$ws.Range("M135").Value = -3643.5
$ws.Range("H137").Value = 1682.3334
$ws.Range("I137").Value = 1124.0834
$ws.Range("J137").Value = 2240.5833
$ws.Range("K137").Value = 3372.2502
$ws.Range("L137").Value = 6721.749899999999
$ws.Range("M137").Value = -822.2501999999999
$ws.Range("N137").Value = -11821.7499
$ws.Range("H138").Value = 2638.04
$ws.Range("I138").Value = 2622.9583
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 7868.874899999999
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = -2728.874899999999
$ws.Range("N138").Value = -19280
$ws.Range("H141").Value = 4663.6665
$ws.Range("I141").Value = 3331.3333
$ws.Range("J141").Value = 5329.8335
$ws.Range("K141").Value = 9993.999899999999
$ws.Range("L141").Value = 15989.5005
$ws.Range("M141").Value = -4813.999899999999
$ws.Range("N141").Value = -26349.5005

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1099.8
$ws.Range("I2").Value = 499.5
$ws.Range("K2").Value = 499.5
$ws.Range("M2").Value = -386.5
$ws.Range("H32").Value = 3734.4126
$ws.Range("I32").Value = 2429.4255
$ws.Range("J32").Value = 7567.8125
$ws.Range("K32").Value = 2429.4255
$ws.Range("L32").Value = 7567.8125
$ws.Range("M32").Value = -2142.4255
$ws.Range("N32").Value = -8141.8125
$ws.Range("H61").Value = 2792.4333
$ws.Range("I61").Value = 1934.3334
$ws.Range("J61").Value = 6224.8335
$ws.Range("K61").Value = 1934.3334
$ws.Range("L61").Value = 6224.8335
$ws.Range("M61").Value = -1722.3334
$ws.Range("N61").Value = -6648.8335
$ws.Range("H116").Value = 1099.8
$ws.Range("I116").Value = 499.5
$ws.Range("K116").Value = 499.5
$ws.Range("M116").Value = 1794.5
$ws.Range("H132").Value = 1232.8096
$ws.Range("I132").Value = 836.3684
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 2509.1052
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = 20.89480000000003
$ws.Range("N132").Value = -20057
$ws.Range("H136").Value = 2792.4333
$ws.Range("I136").Value = 1934.3334
$ws.Range("J136").Value = 6224.8335
$ws.Range("K136").Value = 5803.0002
$ws.Range("L136").Value = 18674.5005
$ws.Range("M136").Value = -3253.0002
$ws.Range("N136").Value = -23774.5005

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1099.8
$ws.Range("I3").Value = 499.5
$ws.Range("K3").Value = 499.5
$ws.Range("M3").Value = -385.5
$ws.Range("H22").Value = 571.5
$ws.Range("J22").Value = 994
$ws.Range("L22").Value = 994
$ws.Range("N22").Value = -1340
$ws.Range("H86").Value = 155733.08
$ws.Range("I86").Value = 1372.5
$ws.Range("J86").Value = 224337.78
$ws.Range("K86").Value = 1372.5
$ws.Range("L86").Value = 224337.78
$ws.Range("M86").Value = -249.5
$ws.Range("N86").Value = -226583.78
$ws.Range("H89").Value = 155733.08
$ws.Range("I89").Value = 1372.5
$ws.Range("J89").Value = 224337.78
$ws.Range("K89").Value = 6862.5
$ws.Range("L89").Value = 1121688.9
$ws.Range("M89").Value = -1246.5
$ws.Range("N89").Value = -1132920.9
$ws.Range("H94").Value = 581.8823
$ws.Range("I94").Value = 505.2857
$ws.Range("K94").Value = 505.2857
$ws.Range("M94").Value = -54.28570000000002
$ws.Range("H107").Value = 1932.8889
$ws.Range("I107").Value = 1799.5
$ws.Range("K107").Value = 1799.5
$ws.Range("M107").Value = 120.5
$ws.Range("H134").Value = 7827.675
$ws.Range("I134").Value = 8523.799999999999
$ws.Range("K134").Value = 25571.4
$ws.Range("M134").Value = -23036.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2795.1667
$ws.Range("I31").Value = 1299.1111
$ws.Range("K31").Value = 1299.1111
$ws.Range("M31").Value = -1004.1111
$ws.Range("H34").Value = 2795.1667
$ws.Range("I34").Value = 1299.1111
$ws.Range("K34").Value = 1299.1111
$ws.Range("M34").Value = -1097.1111
$ws.Range("H58").Value = 1191.2727
$ws.Range("I58").Value = 1026
$ws.Range("K58").Value = 1026
$ws.Range("M58").Value = -823
$ws.Range("H132").Value = 2180
$ws.Range("I132").Value = 1509.8182
$ws.Range("J132").Value = 3101.5
$ws.Range("K132").Value = 4529.4546
$ws.Range("L132").Value = 9304.5
$ws.Range("M132").Value = -1999.4546
$ws.Range("N132").Value = -14364.5
$ws.Range("H134").Value = 843.9375
$ws.Range("I134").Value = 766.86664
$ws.Range("K134").Value = 2300.59992
$ws.Range("M134").Value = 234.4000800000003
$ws.Range("H136").Value = 1191.2727
$ws.Range("I136").Value = 1026
$ws.Range("K136").Value = 3078
$ws.Range("M136").Value = -528

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 188
$ws.Range("I23").Value = 50
$ws.Range("J23").Value = 222.5
$ws.Range("K23").Value = 150
$ws.Range("L23").Value = 667.5
$ws.Range("M23").Value = 85
$ws.Range("N23").Value = -1137.5
$ws.Range("H33").Value = 109.333336
$ws.Range("I33").Value = 121.4
$ws.Range("K33").Value = 728.4000000000001
$ws.Range("M33").Value = -445.4000000000001
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").ClearContents()
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("M121").ClearContents()
$ws.Range("N121").ClearContents()
$ws.Range("H128").Value = 399999.5
$ws.Range("I128").Value = 399999.5
$ws.Range("K128").Value = 1199998.5
$ws.Range("M128").Value = -1195018.5
$ws.Range("H131").Value = 768.77
$ws.Range("I131").Value = 334.55554
$ws.Range("J131").Value = 811.7143
$ws.Range("K131").Value = 1003.66662
$ws.Range("L131").Value = 2435.1429
$ws.Range("M131").Value = 4036.33338
$ws.Range("N131").Value = -12515.1429
$ws.Range("H132").Value = 1466.5834
$ws.Range("I132").Value = 1463.5454
$ws.Range("K132").Value = 13171.9086
$ws.Range("M132").Value = -10641.9086

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5152456
$ws.Range("I11").Value = 7102420
$ws.Range("J11").Value = 2715001
$ws.Range("K11").Value = 7102420
$ws.Range("L11").Value = 2715001
$ws.Range("M11").Value = -7102281
$ws.Range("N11").Value = -2715279
$ws.Range("H12").Value = 5676000.5
$ws.Range("J12").Value = 2380002.8
$ws.Range("L12").Value = 2380002.8
$ws.Range("N12").Value = -2380282.8
$ws.Range("H132").Value = 5541.0557
$ws.Range("I132").Value = 4415.154
$ws.Range("J132").Value = 8468.4
$ws.Range("K132").Value = 13245.462
$ws.Range("L132").Value = 25405.2
$ws.Range("M132").Value = -10715.462
$ws.Range("N132").Value = -30465.2
$ws.Range("H138").Value = 22346.555
$ws.Range("I138").Value = 22346.555
$ws.Range("K138").Value = 22346.555
$ws.Range("M138").Value = -17206.555
$ws.Range("H141").Value = 22476.334
$ws.Range("J141").Value = 22476.334
$ws.Range("L141").Value = 22476.334
$ws.Range("N141").Value = -32836.334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2398
$ws.Range("I16").Value = 4818.5713
$ws.Range("J16").Value = 515.3333
$ws.Range("K16").Value = 4818.5713
$ws.Range("L16").Value = 515.3333
$ws.Range("M16").Value = -4648.5713
$ws.Range("N16").Value = -855.3333
$ws.Range("H20").Value = 14975
$ws.Range("J20").Value = 14975
$ws.Range("L20").Value = 14975
$ws.Range("N20").Value = -15427
$ws.Range("H40").Value = 11699.3
$ws.Range("I40").Value = 6666.3335
$ws.Range("K40").Value = 6666.3335
$ws.Range("M40").Value = -6530.3335
$ws.Range("H100").Value = 1650
$ws.Range("I100").Value = 1650
$ws.Range("K100").Value = 1650
$ws.Range("M100").Value = -1109
$ws.Range("H136").Value = 2694.4333
$ws.Range("I136").Value = 1462.7894
$ws.Range("J136").Value = 4821.8184
$ws.Range("K136").Value = 4388.3682
$ws.Range("L136").Value = 14465.4552
$ws.Range("M136").Value = -1838.3682
$ws.Range("N136").Value = -19565.4552

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 200
$ws.Range("I100").Value = 200
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 400
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = 141
$ws.Range("N100").ClearContents()
$ws.Range("H132").Value = 8022.846
$ws.Range("I132").Value = 2649
$ws.Range("J132").Value = 8999.909
$ws.Range("K132").Value = 7947
$ws.Range("L132").Value = 26999.727
$ws.Range("M132").Value = -5417
$ws.Range("N132").Value = -32059.727
$ws.Range("H136").Value = 2312.1052
$ws.Range("I136").Value = 1748.2222
$ws.Range("J136").Value = 2819.6
$ws.Range("K136").Value = 5244.6666
$ws.Range("L136").Value = 8458.799999999999
$ws.Range("M136").Value = -2694.6666
$ws.Range("N136").Value = -13558.8

Write-Host "Applied 251 cell updates"